# Applies the "Usecases corrected" commit:
#   1. "Use" / " " / "Cases" runs (with spell-check proofErr wrappers around
#      "Use" and "Cases") are merged into "Use" + " Cases" (proofErr kept
#      only around "Use").
#   2. The "Back - Office Auswertung" heading paragraph (and the blank
#      paragraph that followed it) is removed and replaced by a single,
#      style-less paragraph holding Word's "_GoBack" bookmark.
#   3. The "Der MA gibt einen Zeitraum ein. ... bestaetigt hat, wird die
#      gewuenschte Auswertung ..." sentence is re-split into different runs
#      and the grammar-check proofErr markers are moved to wrap "hat, wird"
#      (same visible text, different run/proofing-mark layout).
#   4. The spell-check proofErr wrapper around "Rezeptionist" is removed.
#
# All four are structural (run/proofErr) edits, not plain text substitutions,
# so they are carried out with Range.InsertXML (the supported way to change
# run/paragraph-level markup precisely) rather than Find/Replace.

$d = $word.ActiveDocument

function Get-TextIndex([string]$needle, [int]$start = 0) {
    $full = $d.Content.Text
    return $full.IndexOf($needle, $start)
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>'
$pkgFooter = '</pkg:xmlData></pkg:part></pkg:package>'
$wDoc = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'

function Set-ParagraphXml($range, [string]$bodyXml) {
    $xml = $pkgHeader + $wDoc + '<w:body>' + $bodyXml + '</w:body></w:document>' + $pkgFooter
    $range.InsertXML($xml) | Out-Null
}

# --- 1. "Use Cases" title: drop the proofErr wrapper + separate run around
#        "Cases" and fold the space into the remaining run. -----------------
$titlePara = $d.Paragraphs(1).Range
$newTitle = '<w:p w:rsidR="00912044" w:rsidRDefault="00912044" w:rsidP="00912044"><w:pPr><w:pStyle w:val="Titel"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Use</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Cases</w:t></w:r></w:p>'
Set-ParagraphXml $titlePara $newTitle

# --- 2. Remove the "Back - Office Auswertung" heading paragraph plus the
#        trailing blank paragraph, replacing both with a single paragraph
#        that only contains the "_GoBack" bookmark. ------------------------
$idx = Get-TextIndex "Back - Office Auswertung"
$anchor = $d.Range($idx, $idx + 25)
$headingPara = $anchor.Paragraphs(1)
$blankPara = $headingPara.Next()
$mergedRange = $d.Range($headingPara.Range.Start, $blankPara.Range.End)
$newBack = '<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Set-ParagraphXml $mergedRange $newBack

# --- 3. Re-flow the "MA gibt einen Zeitraum ..." sentence: same visible
#        text, but the grammar-check proofErr pair now wraps
#        "Nachdem der MA ... bestaetigt hat, wird" instead of just "hat". ---
$idx = Get-TextIndex "Der "
$idx = Get-TextIndex "MA gibt einen Zeitraum"
$anchor = $d.Range($idx, $idx + 10)
$sentencePara = $anchor.Paragraphs(1)
$sentenceRange = $sentencePara.Range
$newSentence = '<w:p w:rsidR="00873E72" w:rsidRDefault="00873E72" w:rsidP="00873E72"><w:r><w:t xml:space="preserve">Der </w:t></w:r><w:r w:rsidR="003C669C"><w:t xml:space="preserve">MA gibt einen Zeitraum ein. Au' + [char]0x00DF + 'erdem w' + [char]0x00E4 + 'hlt er Auswahlkriterien und die Art der Sortierung aus. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="003C669C"><w:t xml:space="preserve">Nachdem der MA nun seine Eingaben dem System </w:t></w:r><w:r w:rsidR="00366134"><w:t>best' + [char]0x00E4 + 'tigt hat</w:t></w:r><w:r w:rsidR="003C669C"><w:t>, wird</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="003C669C"><w:t xml:space="preserve"> die gew' + [char]0x00FC + 'nschte Auswertung am Bildschirm ausgegeben.</w:t></w:r></w:p>'
Set-ParagraphXml $sentenceRange $newSentence

# --- 4. Drop the spell-check proofErr wrapper around "Rezeptionist". -------
$idx = Get-TextIndex "Rezeptionist"
$anchor = $d.Range($idx, $idx + 12)
$rezPara = $anchor.Paragraphs(1)
$rezRange = $rezPara.Range
$newRez = '<w:p w:rsidR="00873E72" w:rsidRDefault="00BD4CDE" w:rsidP="00873E72"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Rezeptionist</w:t></w:r></w:p>'
Set-ParagraphXml $rezRange $newRez

Write-Host "edit.ps1 complete"
